$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1228
$ws.Range("F5").Value = 1451
$ws.Range("F6").Value = 1737
$ws.Range("F7").Value = 6267
$ws.Range("F8").Value = 131
$ws.Range("F9").Value = 1861
$ws.Range("F10").Value = 494
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = 42
$ws.Range("F16").Value = 7091
$ws.Range("F17").Value = 133
$ws.Range("F19").Value = 177
$ws.Range("F21").Value = 1723
$ws.Range("F28").Value = 1652
$ws.Range("F30").Value = 337
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F25").Value = 2
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9538
$ws.Range("F3").Value = 2263
$ws.Range("F5").Value = 255
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9538
$ws.Range("F3").Value = 2263
$ws.Range("F5").Value = 1228
$ws.Range("F10").Value = 1451
$ws.Range("F11").Value = 255
$ws.Range("F12").Value = 1737
$ws.Range("F13").Value = 6267
$ws.Range("F14").Value = 131
$ws.Range("F15").Value = 1861
$ws.Range("F18").Value = 494
$ws.Range("F21").Value = 4
$ws.Range("F23").Value = 42
$ws.Range("F24").Value = 7091
$ws.Range("F25").Value = 133
$ws.Range("F27").Value = 177
$ws.Range("F29").Value = 1723
$ws.Range("F34").Value = 1652
$ws.Range("F36").Value = 337
$ws.Range("F49").Value = 2
